$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet "NextBus1" (sheet1.xml) - value-only updates, no new rows
# ------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("NextBus1")

$ws1.Range("F2").Value = 45688.60628472222
$ws1.Range("O2").Value = 10

$ws1.Range("F3").Value = 45688.61091435186
$ws1.Range("L3").Value = "DD"
$ws1.Range("O3").Value = 17

$ws1.Range("F4").Value = 45688.61371527778
$ws1.Range("O4").Value = 21

$ws1.Range("F5").Value = 45688.61773148148
$ws1.Range("O5").Value = 26

$ws1.Range("F6").Value = 45688.62018518519
$ws1.Range("O6").Value = 30

$ws1.Range("F7").Value = 45688.61545138889
$ws1.Range("L7").Value = "SD"
$ws1.Range("O7").Value = 23

$ws1.Range("F8").Value = 45688.60774305555
$ws1.Range("O8").Value = 12

# ------------------------------------------------------------------
# Sheet "NextBus2" (sheet2.xml) - value updates on rows 2-6, plus a
# brand-new row 7
# ------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("NextBus2")

$ws2.Range("F2").Value = 45688.61583333334
$ws2.Range("O2").Value = 24

$ws2.Range("F3").Value = 45688.61922453704
$ws2.Range("O3").Value = 29

$ws2.Range("F4").Value = 45688.62368055555
$ws2.Range("O4").Value = 35

$ws2.Range("F5").Value = 45688.62696759259
$ws2.Range("O5").Value = 40

$ws2.Range("B6").Value = 151
$ws2.Range("F6").Value = 45688.62586805555
$ws2.Range("K6").Value = 16009
$ws2.Range("O6").Value = 38

# New row 7, mirroring the row-6 layout/format first, then filling values
$ws2.Range("A6:O6").Copy()
$ws2.Range("A7:O7").PasteSpecial(-4122)
$ws2.Range("A7").Value = "NextBus3"
$ws2.Range("B7").Value = 74
$ws2.Range("C7").Value = 64009
$ws2.Range("D7").Value = "Hougang Ctrl Int"
$ws2.Range("E7").Value = "SBST"
$ws2.Range("F7").Value = 45688.61568287037
$ws2.Range("G7").Value = 64009
$ws2.Range("H7").Value = "WAB"
$ws2.Range("I7").Value = "SEA"
$ws2.Range("J7").Value = 1
$ws2.Range("K7").Value = 11379
$ws2.Range("L7").Value = "DD"
$ws2.Range("M7").Value = 12101
$ws2.Range("N7").Value = "Ngee Ann Poly"
$ws2.Range("O7").Value = 24

# ------------------------------------------------------------------
# Sheet "NextBus3" (sheet3.xml) - identical data edits to NextBus2
# ------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("NextBus3")

$ws3.Range("F2").Value = 45688.61583333334
$ws3.Range("O2").Value = 24

$ws3.Range("F3").Value = 45688.61922453704
$ws3.Range("O3").Value = 29

$ws3.Range("F4").Value = 45688.62368055555
$ws3.Range("O4").Value = 35

$ws3.Range("F5").Value = 45688.62696759259
$ws3.Range("O5").Value = 40

$ws3.Range("B6").Value = 151
$ws3.Range("F6").Value = 45688.62586805555
$ws3.Range("K6").Value = 16009
$ws3.Range("O6").Value = 38

$ws3.Range("A6:O6").Copy()
$ws3.Range("A7:O7").PasteSpecial(-4122)
$ws3.Range("A7").Value = "NextBus3"
$ws3.Range("B7").Value = 74
$ws3.Range("C7").Value = 64009
$ws3.Range("D7").Value = "Hougang Ctrl Int"
$ws3.Range("E7").Value = "SBST"
$ws3.Range("F7").Value = 45688.61568287037
$ws3.Range("G7").Value = 64009
$ws3.Range("H7").Value = "WAB"
$ws3.Range("I7").Value = "SEA"
$ws3.Range("J7").Value = 1
$ws3.Range("K7").Value = 11379
$ws3.Range("L7").Value = "DD"
$ws3.Range("M7").Value = 12101
$ws3.Range("N7").Value = "Ngee Ann Poly"
$ws3.Range("O7").Value = 24
